$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Match the recorded window-position tweak from the original edit (xWindow 3720 -> 4650).
# (Best-effort; window chrome position is cosmetic and may not round-trip to the saved
# bookViews/workbookView xWindow attribute in this environment.)
try {
  $excel.ActiveWindow.Left = 4650
} catch {
}

# Clear out the old data range entirely first.
$ws.Range("A1:D4").ClearContents()

$ws.Range("A1").Value = "# hash marks indicate file comments and are skipped."
$ws.Range("A2").Value = "# '^' indicate attribute names."
$ws.Range("A3").Value = "# Attribute types are parsed and determined to be floating-point, integer, arbitrary or [enumeration]"

$ws.Range("A5").Value = "#Attribute Names:"
$ws.Range("B5").Value = "^temperature"
$ws.Range("C5").Value = "^capacity"
$ws.Range("D5").Value = "^color"
$ws.Range("E5").Value = "^comment"

$ws.Range("A6").Value = "#Sample 1"
$ws.Range("B6").Value = 25.5
$ws.Range("C6").Value = 25
$ws.Range("D6").Value = "comment 1"
$ws.Range("E6").Value = "[red]"

$ws.Range("A7").Value = "#Sample 2"
$ws.Range("B7").Value = 23.5
$ws.Range("C7").Value = 23
$ws.Range("D7").Value = "comment 2"
$ws.Range("E7").Value = "[green]"

$ws.Range("A8").Value = "#Sample 3"
$ws.Range("B8").Value = 10.5
$ws.Range("C8").Value = 10
$ws.Range("D8").Value = "comment 3"
$ws.Range("E8").Value = "[blue]"

$ws.Range("A9").Value = "#Sample 4"
$ws.Range("B9").Value = 5.5
$ws.Range("C9").Value = 5
$ws.Range("D9").Value = "comment 4"
$ws.Range("E9").Value = "[yellow]"

$ws.Range("D10").Select()
